$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("isa_template")

# Row 13 holds the Tags list. The "Computational Analysis" tag (which had no
# ontology accession/source) is dropped, the remaining tags shift one column
# to the left, and a new "Computation" tag (with ontology IDs) is appended.
$ws.Range("D13").Value2 = "mandatory"
$ws.Range("E13").Value2 = "Data Processing"
$ws.Range("F13").Value2 = "Computation"

# Row 14 holds the Tags Term Accession Number, aligned under row 13.
$ws.Range("E14").Value2 = "http://purl.obolibrary.org/obo/NCIT_C47925"
$ws.Range("F14").Value2 = "http://purl.obolibrary.org/obo/NCIT_C61298"

# Row 15 holds the Tags Term Source REF, aligned under row 13.
$ws.Range("E15").Value2 = "NCIT"
$ws.Range("F15").Value2 = "NCIT"
